$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.77062261171358
$ws.Range("C2").Value = 8.28909245029398
$ws.Range("D2").Value = 11.02070741106156
$ws.Range("F2").Value = 31.27156940681103
$ws.Range("G2").Value = 3.630599936758346
$ws.Range("J2").Value = 11.31532976177555
$ws.Range("N2").Value = 16.72606001585262
$ws.Range("O2").Value = 22.50128920119138
$ws.Range("B3").Value = 14.22076246652055
$ws.Range("C3").Value = 7.793167420331984
$ws.Range("D3").Value = 10.94845887686247
$ws.Range("F3").Value = 31.22979351790346
$ws.Range("G3").Value = 3.632935291580253
$ws.Range("J3").Value = 11.29702547725283
$ws.Range("N3").Value = 16.77463508004433
$ws.Range("O3").Value = 22.53699519947344
$ws.Range("B4").Value = 13.87404117160827
$ws.Range("C4").Value = 7.471839705447739
$ws.Range("D4").Value = 10.90590904532199
$ws.Range("F4").Value = 31.21367588684572
$ws.Range("G4").Value = 3.634444977134963
$ws.Range("J4").Value = 11.28812540311812
$ws.Range("N4").Value = 16.80634236683635
$ws.Range("O4").Value = 22.56538612541161
$ws.Range("B5").Value = 13.73067638935665
$ws.Range("C5").Value = 7.336731380094379
$ws.Range("D5").Value = 10.8890392604194
$ws.Range("F5").Value = 31.20950729660418
$ws.Range("G5").Value = 3.635079301762444
$ws.Range("J5").Value = 11.28508895291286
$ws.Range("N5").Value = 16.81973741789058
$ws.Range("O5").Value = 22.57857616516829
$ws.Range("B6").Value = 13.7067523331536
$ws.Range("C6").Value = 7.314047489722779
$ws.Range("D6").Value = 10.88626680967515
$ws.Range("F6").Value = 31.20896004914019
$ws.Range("G6").Value = 3.635185787265117
$ws.Range("J6").Value = 11.2846204691244
$ws.Range("N6").Value = 16.82199031359592
$ws.Range("O6").Value = 22.58086405383218
$ws.Range("B7").Value = 13.87211579530327
$ws.Range("C7").Value = 7.470034349614972
$ws.Range("D7").Value = 10.90567961383011
$ws.Range("F7").Value = 31.21360995115057
$ws.Range("G7").Value = 3.634453454377544
$ws.Range("J7").Value = 11.28808205930267
$ws.Range("N7").Value = 16.8065210965699
$ws.Range("O7").Value = 22.56555745770702
$ws.Range("B8").Value = 14.58304359547821
$ws.Range("C8").Value = 8.12162532492375
$ws.Range("D8").Value = 10.99542897899306
$ws.Range("F8").Value = 31.2551886396931
$ws.Range("G8").Value = 3.63138947936532
$ws.Range("J8").Value = 11.30853452648393
$ws.Range("N8").Value = 16.74241858385625
$ws.Range("O8").Value = 22.51225509723146
$ws.Range("B9").Value = 15.89652648019193
$ws.Range("C9").Value = 9.263677955568763
$ws.Range("D9").Value = 11.18513980608671
$ws.Range("F9").Value = 31.41215400470238
$ws.Range("G9").Value = 3.625979347569381
$ws.Range("J9").Value = 11.36707972421557
$ws.Range("N9").Value = 16.63160971626749
$ws.Range("O9").Value = 22.45927026538271
$ws.Range("B10").Value = 16.80260274048421
$ws.Range("C10").Value = 10.03797756844444
$ws.Range("D10").Value = 11.3319571137071
$ws.Range("F10").Value = 31.57302327772127
$ws.Range("G10").Value = 3.622365249793971
$ws.Range("J10").Value = 11.42114507562969
$ws.Range("N10").Value = 16.5592292481991
$ws.Range("O10").Value = 22.45202492212481
$ws.Range("B11").Value = 17.20032541560146
$ws.Range("C11").Value = 10.37563825539757
$ws.Range("D11").Value = 11.40015560932114
$ws.Range("F11").Value = 31.65595262716799
$ws.Range("G11").Value = 3.620798575083069
$ws.Range("J11").Value = 11.44809104021795
$ws.Range("N11").Value = 16.528252239397
$ws.Range("O11").Value = 22.4556475017766
$ws.Range("B12").Value = 17.34873747353617
$ws.Range("C12").Value = 10.50049100675318
$ws.Range("D12").Value = 11.42616538529199
$ws.Range("F12").Value = 31.68874164064083
$ws.Range("G12").Value = 3.620216380571121
$ws.Range("J12").Value = 11.45862793703306
$ws.Range("N12").Value = 16.51680162510798
$ws.Range("O12").Value = 22.45801587671632
$ws.Range("B13").Value = 17.31687370445016
$ws.Range("C13").Value = 10.47373550671774
$ws.Range("D13").Value = 11.42055580467926
$ws.Range("F13").Value = 31.68161860770888
$ws.Range("G13").Value = 3.620341275104829
$ws.Range("J13").Value = 11.45634390106667
$ws.Range("N13").Value = 16.51925528885777
$ws.Range("O13").Value = 22.45746146399746
$ws.Range("B14").Value = 17.21258001510228
$ws.Range("C14").Value = 10.38597026482095
$ws.Range("D14").Value = 11.40229185795929
$ws.Range("F14").Value = 31.65862253475935
$ws.Range("G14").Value = 3.620750456039642
$ws.Range("J14").Value = 11.44895127300221
$ws.Range("N14").Value = 16.52730458801744
$ws.Range("O14").Value = 22.4558223698615
$ws.Range("B15").Value = 17.148407780252
$ws.Range("C15").Value = 10.33181961226104
$ws.Range("D15").Value = 11.39112813935358
$ws.Range("F15").Value = 31.64471667664368
$ws.Range("G15").Value = 3.621002531249319
$ws.Range("J15").Value = 11.4444662888603
$ws.Range("N15").Value = 16.53227142516748
$ws.Range("O15").Value = 22.45494819701654
$ws.Range("B16").Value = 16.77630927716089
$ws.Range("C16").Value = 10.01548852771207
$ws.Range("D16").Value = 11.32752700366384
$ws.Range("F16").Value = 31.5677984755446
$ws.Range("G16").Value = 3.622469188090947
$ws.Range("J16").Value = 11.41943098495924
$ws.Range("N16").Value = 16.56129284100247
$ws.Range("O16").Value = 22.4519275630465
$ws.Range("B17").Value = 16.54424724709641
$ws.Range("C17").Value = 9.81605251191716
$ws.Range("D17").Value = 11.28885808470129
$ws.Range("F17").Value = 31.52309769119038
$ws.Range("G17").Value = 3.623388716028431
$ws.Range("J17").Value = 11.40467140612282
$ws.Range("N17").Value = 16.57959535575567
$ws.Range("O17").Value = 22.45184800368563
$ws.Range("B18").Value = 16.40941664284531
$ws.Range("C18").Value = 9.699362477591695
$ws.Range("D18").Value = 11.26675063403003
$ws.Range("F18").Value = 31.49830575208768
$ws.Range("G18").Value = 3.623924892560912
$ws.Range("J18").Value = 11.39640374100581
$ws.Range("N18").Value = 16.59030598239102
$ws.Range("O18").Value = 22.45245340205987
$ws.Range("B19").Value = 16.36353651278694
$ws.Range("C19").Value = 9.659512881655075
$ws.Range("D19").Value = 11.25928898095512
$ws.Range("F19").Value = 31.49006988908832
$ws.Range("G19").Value = 3.624107686292894
$ws.Range("J19").Value = 11.39364266843591
$ws.Range("N19").Value = 16.59396395199628
$ws.Range("O19").Value = 22.45277014990164
$ws.Range("B20").Value = 16.56909175454516
$ws.Range("C20").Value = 9.837487628277247
$ws.Range("D20").Value = 11.29296073218225
$ws.Range("F20").Value = 31.52776118827672
$ws.Range("G20").Value = 3.623290076813316
$ws.Range("J20").Value = 11.40621968348493
$ws.Range("N20").Value = 16.5776280341582
$ws.Range("O20").Value = 22.45178906512852
$ws.Range("B21").Value = 17.24327408689077
$ws.Range("C21").Value = 10.41183069601909
$ws.Range("D21").Value = 11.40765156254887
$ws.Range("F21").Value = 31.66533958033622
$ws.Range("G21").Value = 3.620629969673943
$ws.Range("J21").Value = 11.45111367098324
$ws.Range("N21").Value = 16.5249327271501
$ws.Range("O21").Value = 22.45627675547275
$ws.Range("B22").Value = 17.67104196553528
$ws.Range("C22").Value = 10.76964633657051
$ws.Range("D22").Value = 11.48367528887198
$ws.Range("F22").Value = 31.76332130787371
$ws.Range("G22").Value = 3.618955940539192
$ws.Range("J22").Value = 11.48239306388591
$ws.Range("N22").Value = 16.49212338704491
$ws.Range("O22").Value = 22.46501871394139
$ws.Range("B23").Value = 17.44394497333373
$ws.Range("C23").Value = 10.58027540397405
$ws.Range("D23").Value = 11.44300852674462
$ws.Range("F23").Value = 31.71029458747441
$ws.Range("G23").Value = 3.61984351837557
$ws.Range("J23").Value = 11.46552307417142
$ws.Range("N23").Value = 16.50948538424607
$ws.Range("O23").Value = 22.45982111450323
$ws.Range("B24").Value = 16.55786395113678
$ws.Range("C24").Value = 9.827803141125045
$ws.Range("D24").Value = 11.29110553815974
$ws.Range("F24").Value = 31.52564999428052
$ws.Range("G24").Value = 3.623334648149184
$ws.Range("J24").Value = 11.40551902844699
$ws.Range("N24").Value = 16.57851687374051
$ws.Range("O24").Value = 22.45181368307131
$ws.Range("B25").Value = 15.55091524400509
$ws.Range("C25").Value = 8.969508439778641
$ws.Range("D25").Value = 11.13244299755601
$ws.Range("F25").Value = 31.36165079438571
$ws.Range("G25").Value = 3.627379295265041
$ws.Range("J25").Value = 11.34928533489924
$ws.Range("N25").Value = 16.65999702044128
$ws.Range("O25").Value = 22.46805599524075
